$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 69, shifting rows 69:145 down to 70:146
$ws.Rows("69:69").Insert()

# Populate the newly inserted row 69 with the new data entry
$ws.Range("A69").Value = 11
$ws.Range("B69").Value = "Vega Monumental Concepción"
$ws.Range("C69").Value = "Bíobío"
$ws.Range("D69").Value = 45225
$ws.Range("E69").Value = 8
$ws.Range("F69").Value = 100112012
$ws.Range("G69").Value = "Espinaca"
$ws.Range("H69").Value = "Sin especificar"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 80
$ws.Range("K69").Value = 10000
$ws.Range("L69").Value = 10000
$ws.Range("M69").Value = 10000
$ws.Range("N69").Value = "`$/cuna 10 kilos"
$ws.Range("O69").Value = "Región Metropolitana"
$ws.Range("P69").Value = 1000
$ws.Range("Q69").Value = 10
$ws.Range("R69").Value = "Hortaliza"
